$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their original text storage type
$ws.Range("D2:E51").NumberFormat = "@"

# Special case: D31 contains a subscript-3 unicode character (U+2083)
$sub3 = [string][char]8323
$ws.Range("D31").Value = [string]::Concat("0.0", $sub3, "0739")

$ws.Range("D2").Value = "57.325.39"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "2.329.40"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "540.26"
$ws.Range("E5").Value = "  +4.84%  "
$ws.Range("D6").Value = "135.42"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").Value = "0.993"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("D8").Value = "0.535"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "2.363.81"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "0.353"
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("D14").Value = "23.71"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "2.752.07"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "57.416.93"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "2.355.62"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "337.44"
$ws.Range("E19").Value = "  +3.14%  "
$ws.Range("D20").Value = "10.50"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "61.56"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").Value = "8.47"
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("D27").Value = "0.991"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("E28").Value = "  +4.74%  "
$ws.Range("D29").Value = "173.92"
$ws.Range("E29").Value = "  +3.20%  "
$ws.Range("D30").Value = "1.75"
$ws.Range("E30").Value = "  +4.13%  "
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").Value = "6.14"
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").Value = "18.49"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +11.76%  "
$ws.Range("D36").Value = "0.991"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").Value = "4.10"
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("D40").Value = "39.34"
$ws.Range("D41").Value = "149.22"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("D43").Value = "3.63"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").Value = "285.99"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("D45").Value = "0.0928"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "18.92"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.562"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("D50").Value = "17.53"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("E51").Value = "  +7.24%  "
